$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("G2").Value = 201.4397426666667
$ws.Range("H2").Value = 604.3192280000001
$ws.Range("I2").Value = 0.4833500233086392
$ws.Range("J2").Value = 0.4833500233086393
$ws.Range("M2").Value = 0.1137376666666667
$ws.Range("N2").Value = 0.341213
$ws.Range("O2").Value = 0.003048021899328029
$ws.Range("P2").Value = 0.003048021899328029
$ws.Range("Q2").Value = 22.91128630484045
$ws.Range("R2").Value = 206.201576743564
$ws.Range("S2").Value = 0.001473261456085446
$ws.Range("T2").Value = 0.001473261456085446
# Row 3
$ws.Range("G3").Value = 201.4397426666667
$ws.Range("H3").Value = 604.3192280000001
$ws.Range("I3").Value = 0.4833500233086392
$ws.Range("J3").Value = 0.4833500233086393
$ws.Range("O3").Value = 0.0144044366216848
$ws.Range("P3").Value = 0.0144044366216848
$ws.Range("Q3").Value = 108.2748688820467
$ws.Range("R3").Value = 974.4738199384201
$ws.Range("S3").Value = 0.006962384776839167
$ws.Range("T3").Value = 0.006962384776839168
# Row 4
$ws.Range("G4").Value = 201.4397426666667
$ws.Range("H4").Value = 604.3192280000001
$ws.Range("I4").Value = 0.4833500233086392
$ws.Range("J4").Value = 0.4833500233086393
$ws.Range("K4").Value = 1
$ws.Range("L4").Value = 0.3333333333333333
$ws.Range("M4").Value = 0.05518366666666667
$ws.Range("N4").Value = 0.165551
$ws.Range("O4").Value = 0.001478850669393178
$ws.Range("P4").Value = 0.001478850669393178
$ws.Range("Q4").Value = 11.11618361273645
$ws.Range("R4").Value = 100.045652514628
$ws.Range("S4").Value = 0.0007148025055211895
$ws.Range("T4").Value = 0.0007148025055211895
# Row 5
$ws.Range("G5").Value = 201.4397426666667
$ws.Range("H5").Value = 604.3192280000001
$ws.Range("I5").Value = 0.4833500233086392
$ws.Range("J5").Value = 0.4833500233086393
$ws.Range("M5").Value = 36.608813
$ws.Range("N5").Value = 109.826439
$ws.Range("O5").Value = 0.981068690809594
$ws.Range("P5").Value = 0.9810686908095939
$ws.Range("Q5").Value = 7374.469870052123
$ws.Range("R5").Value = 66370.2288304691
$ws.Range("S5").Value = 0.4741995745701935
$ws.Range("T5").Value = 0.4741995745701935
# Row 6
$ws.Range("I6").Value = 0.1569674599353791
$ws.Range("J6").Value = 0.1569674599353792
$ws.Range("M6").Value = 0.1137376666666667
$ws.Range("N6").Value = 0.341213
$ws.Range("O6").Value = 0.003048021899328029
$ws.Range("P6").Value = 0.003048021899328029
$ws.Range("Q6").Value = 7.440418416669112
$ws.Range("R6").Value = 66.96376575002199
$ws.Range("S6").Value = 0.0004784402553649307
$ws.Range("T6").Value = 0.0004784402553649308
# Row 7
$ws.Range("I7").Value = 0.1569674599353791
$ws.Range("J7").Value = 0.1569674599353792
$ws.Range("O7").Value = 0.0144044366216848
$ws.Range("P7").Value = 0.0144044366216848
$ws.Range("S7").Value = 0.002261027828306018
$ws.Range("T7").Value = 0.002261027828306018
# Row 8
$ws.Range("I8").Value = 0.1569674599353791
$ws.Range("J8").Value = 0.1569674599353792
$ws.Range("K8").Value = 1
$ws.Range("L8").Value = 0.3333333333333333
$ws.Range("M8").Value = 0.05518366666666667
$ws.Range("N8").Value = 0.165551
$ws.Range("O8").Value = 0.001478850669393178
$ws.Range("P8").Value = 0.001478850669393178
$ws.Range("Q8").Value = 3.609970045977111
$ws.Range("R8").Value = 32.489730413794
$ws.Range("S8").Value = 0.0002321314331983823
$ws.Range("T8").Value = 0.0002321314331983824
# Row 9
$ws.Range("I9").Value = 0.1569674599353791
$ws.Range("J9").Value = 0.1569674599353792
$ws.Range("M9").Value = 36.608813
$ws.Range("N9").Value = 109.826439
$ws.Range("O9").Value = 0.981068690809594
$ws.Range("P9").Value = 0.9810686908095939
$ws.Range("Q9").Value = 2394.852070034808
$ws.Range("R9").Value = 21553.66863031327
$ws.Range("S9").Value = 0.1539958604185098
$ws.Range("T9").Value = 0.1539958604185098
# Row 10
$ws.Range("G10").Value = 60.43484133333334
$ws.Range("H10").Value = 181.304524
$ws.Range("I10").Value = 0.1450120099461104
$ws.Range("J10").Value = 0.1450120099461104
$ws.Range("M10").Value = 0.1137376666666667
$ws.Range("N10").Value = 0.341213
$ws.Range("O10").Value = 0.003048021899328029
$ws.Range("P10").Value = 0.003048021899328029
$ws.Range("Q10").Value = 6.873717838623556
$ws.Range("R10").Value = 61.863460547612
$ws.Range("S10").Value = 0.0004419997819813183
$ws.Range("T10").Value = 0.0004419997819813185
# Row 11
$ws.Range("G11").Value = 60.43484133333334
$ws.Range("H11").Value = 181.304524
$ws.Range("I11").Value = 0.1450120099461104
$ws.Range("J11").Value = 0.1450120099461104
$ws.Range("O11").Value = 0.0144044366216848
$ws.Range("P11").Value = 0.0144044366216848
$ws.Range("Q11").Value = 32.48402939087334
$ws.Range("R11").Value = 292.35626451786
$ws.Range("S11").Value = 0.002088816306651873
$ws.Range("T11").Value = 0.002088816306651873
# Row 12
$ws.Range("G12").Value = 60.43484133333334
$ws.Range("H12").Value = 181.304524
$ws.Range("I12").Value = 0.1450120099461104
$ws.Range("J12").Value = 0.1450120099461104
$ws.Range("K12").Value = 1
$ws.Range("L12").Value = 0.3333333333333333
$ws.Range("M12").Value = 0.05518366666666667
$ws.Range("N12").Value = 0.165551
$ws.Range("O12").Value = 0.001478850669393178
$ws.Range("P12").Value = 0.001478850669393178
$ws.Range("Q12").Value = 3.335016139191556
$ws.Range("R12").Value = 30.015145252724
$ws.Range("S12").Value = 0.0002144511079788555
$ws.Range("T12").Value = 0.0002144511079788556
# Row 13
$ws.Range("G13").Value = 60.43484133333334
$ws.Range("H13").Value = 181.304524
$ws.Range("I13").Value = 0.1450120099461104
$ws.Range("J13").Value = 0.1450120099461104
$ws.Range("M13").Value = 36.608813
$ws.Range("N13").Value = 109.826439
$ws.Range("O13").Value = 0.981068690809594
$ws.Range("P13").Value = 0.9810686908095939
$ws.Range("Q13").Value = 2212.447805056671
$ws.Range("R13").Value = 19912.03024551004
$ws.Range("S13").Value = 0.1422667427494983
$ws.Range("T13").Value = 0.1422667427494983
# Row 14
$ws.Range("G14").Value = 89.46554166666668
$ws.Range("H14").Value = 268.396625
$ws.Range("I14").Value = 0.2146705068098712
$ws.Range("J14").Value = 0.2146705068098712
$ws.Range("M14").Value = 0.1137376666666667
$ws.Range("N14").Value = 0.341213
$ws.Range("O14").Value = 0.003048021899328029
$ws.Range("P14").Value = 0.003048021899328029
$ws.Range("Q14").Value = 10.17560195623611
$ws.Range("R14").Value = 91.58041760612501
$ws.Range("S14").Value = 0.0006543204058963342
$ws.Range("T14").Value = 0.0006543204058963343
# Row 15
$ws.Range("G15").Value = 89.46554166666668
$ws.Range("H15").Value = 268.396625
$ws.Range("I15").Value = 0.2146705068098712
$ws.Range("J15").Value = 0.2146705068098712
$ws.Range("O15").Value = 0.0144044366216848
$ws.Range("P15").Value = 0.0144044366216848
$ws.Range("Q15").Value = 48.08817597354167
$ws.Range("R15").Value = 432.793583761875
$ws.Range("S15").Value = 0.003092207709887746
$ws.Range("T15").Value = 0.003092207709887746
# Row 16
$ws.Range("G16").Value = 89.46554166666668
$ws.Range("H16").Value = 268.396625
$ws.Range("I16").Value = 0.2146705068098712
$ws.Range("J16").Value = 0.2146705068098712
$ws.Range("K16").Value = 1
$ws.Range("L16").Value = 0.3333333333333333
$ws.Range("M16").Value = 0.05518366666666667
$ws.Range("N16").Value = 0.165551
$ws.Range("O16").Value = 0.001478850669393178
$ws.Range("P16").Value = 0.001478850669393178
$ws.Range("Q16").Value = 4.937036629486112
$ws.Range("R16").Value = 44.43332966537501
$ws.Range("S16").Value = 0.0003174656226947509
$ws.Range("T16").Value = 0.0003174656226947509
# Row 17
$ws.Range("G17").Value = 89.46554166666668
$ws.Range("H17").Value = 268.396625
$ws.Range("I17").Value = 0.2146705068098712
$ws.Range("J17").Value = 0.2146705068098712
$ws.Range("M17").Value = 36.608813
$ws.Range("N17").Value = 109.826439
$ws.Range("O17").Value = 0.981068690809594
$ws.Range("P17").Value = 0.9810686908095939
$ws.Range("Q17").Value = 3275.227284818709
$ws.Range("R17").Value = 29477.04556336838
$ws.Range("S17").Value = 0.2106065130713924
$ws.Range("T17").Value = 0.2106065130713924
